$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string table gets "line7"/"line8" spliced in right after "line6"
# (i.e. before "extr1"), but none of the existing rows' string references are
# renumbered to compensate. The net, visible effect is that every row from
# the old "extr1" row onward is relabelled two slots down the name list
# (extr1->line7, extr2->line8, extr3->extr1, ... extr8->extr6), while two
# brand-new rows are appended carrying the labels extr7/extr8. Combined with
# the numeric tweaks in the diff, the final row-by-row content is:

# Row 8: extr1 -> line7; C 5->14, D 12->11, E False->True
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9: extr2 -> line8; C 5->16 (D, E unchanged)
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16

# Row 10: extr3 -> extr1; C 10->5, D 11->12 (E unchanged, stays True)
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# Row 11: extr4 -> extr2; C 7->5, D 8->9 (E unchanged, stays True)
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# Row 12: extr5 -> extr3; C 9->10, E False->True (D unchanged)
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true

# Row 13: extr6 -> extr4; D 11->8, E False->True (C unchanged)
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 14: extr7 -> extr5; C 5->9, D 7->11, E True->False
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# Row 15: extr8 -> extr6; C 8->7, D 5->11, E unchanged (stays False)
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- Append two new rows (16: extr7, 17: extr8), matching column A's style ---

# Row 16
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# Row 17
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

$excel.CutCopyMode = 0
